$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the missing X6/Y6 values that complete row 6 ---
$ws.Range("X6").Value = -0.35999999999999943
$ws.Range("Y6").Value = "Down"

# --- Append new row 7 with a full set of scan values ---
$ws.Range("A7").Value = 42648.885243055556
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"

$ws.Range("B7").Value = -8
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 20869
$ws.Range("F7").Value = 3280
$ws.Range("G7").Value = 57
$ws.Range("H7").Value = 42
$ws.Range("I7").Value = 69
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 26054
$ws.Range("L7").Value = 385
$ws.Range("M7").Value = 281
$ws.Range("N7").Value = 73
$ws.Range("O7").Value = 32
$ws.Range("P7").Value = "Noun"
$ws.Range("Q7").Value = 52.976913006825477
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = -0.094399999999999998
$ws.Range("S7").NumberFormat = "0.00%"
$ws.Range("T7").Value = -0.025700000000000001
$ws.Range("T7").NumberFormat = "0.00%"
$ws.Range("U7").Value = 6.62
$ws.Range("V7").Value = 1.88
$ws.Range("W7").Value = -2
